$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.039.22"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "1.680.35"
$ws.Range("E3").Value = "  +0.85%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("E6").Value = "  -2.48%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.254"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0624"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("E11").Value = "  -0.75%  "

$ws.Range("D12").Value = "1.917.47"
$ws.Range("E12").Value = "  +0.89%  "

$ws.Range("D13").Value = "1.694.89"
$ws.Range("E13").Value = "  +1.78%  "

$ws.Range("E14").Value = "  +0.82%  "

$ws.Range("E15").Value = "  +1.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("D17").Value = "27.035.92"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "235.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("E20").Value = "  +0.81%  "

$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("E22").Value = "  +3.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.76%  "

$ws.Range("E24").Value = "  -3.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.67%  "

$ws.Range("E28").Value = "  -2.55%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +0.21%  "

$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").Value = "1.540.58"
$ws.Range("E33").Value = "  +5.48%  "

$ws.Range("E34").Value = "  +1.90%  "

$ws.Range("E35").Value = "  +5.56%  "

$ws.Range("E36").Value = "  -0.83%  "

$ws.Range("E37").Value = "  +1.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.916"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.14%  "

$ws.Range("E39").Value = "  +3.02%  "

$ws.Range("E40").Value = "  +6.29%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.42%  "

$ws.Range("E44").Value = "  -0.38%  "

$ws.Range("D45").Value = "1.823.00"
$ws.Range("E45").Value = "  +0.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.780"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("E48").Value = "  +0.25%  "

$ws.Range("E49").Value = "  +2.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.36%  "

$ws.Range("E51").Value = "  +0.18%  "
